$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data row 17 (all numeric values)
$ws.Range("A17").Value = 20250725
$ws.Range("B17").Value = 1
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 4
$ws.Range("F17").Value = 5
$ws.Range("G17").Value = 6

# New row 18 (string trial-id lists)
$ws.Range("E18").Value = "10,13,14,17,18"
$ws.Range("F18").Value = "4,5,9,12,17,19"
$ws.Range("G18").Value = "5,9,13,14,18,19"

# Update selection to match the diff (G18 active cell)
$ws.Range("G18").Select()
